$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 33 (copying row 32's formatting) for the new item ---
$ws.Rows.Item(32).Copy()
$ws.Rows.Item(33).Insert()

# Fix up the bottom border that Insert() drops, so the engine re-uses the
# existing style indices (7,8,9,10,11,12) instead of minting new ones.
$rng = $ws.Range("A33:Q33")
$bd = $rng.Borders.Item(9)
$bd.Color = 13882323
$bd.LineStyle = 1

# Row heights: new data row matches the old "totals" row height (24.75),
# and the totals row (now at 34) becomes 25.5.
$ws.Rows.Item(33).RowHeight = 24.75
$ws.Rows.Item(34).RowHeight = 25.5

# --- 2. Populate the new row 33 with the new item's data ---
# (helper: remember each cell's number format, force text entry so
#  numeric-looking strings like "35.0000" are kept as shared-string text
#  instead of being auto-converted to numbers, then restore the format)
function Set-TextValue($cell, $text) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.NumberFormat = $fmt
}

$ws.Range("A33").Value2 = 27
$ws.Range("B33").Value2 = ""
Set-TextValue $ws.Range("C33") "معجون اسنان كلوس اب وسط"
$ws.Range("D33").Value2 = ""
$ws.Range("E33").Value2 = ""
$ws.Range("F33").Value2 = ""
$ws.Range("G33").Value2 = ""
Set-TextValue $ws.Range("H33") "7:0"
$ws.Range("I33").Value2 = ""
$ws.Range("J33").Value2 = ""
$ws.Range("K33").Value2 = ""
Set-TextValue $ws.Range("L33") "0"
$ws.Range("M33").Value2 = ""
Set-TextValue $ws.Range("N33") "35.00"
$ws.Range("O33").Value2 = ""
Set-TextValue $ws.Range("P33") "35.0000"
Set-TextValue $ws.Range("Q33") "1:0"

# --- 3. Add merge cells for the new row 33 (mirrors row 32's merges) ---
$ws.Range("A33:B33").Merge()
$ws.Range("C33:G33").Merge()
$ws.Range("H33:K33").Merge()
$ws.Range("L33:M33").Merge()
$ws.Range("N33:O33").Merge()

# --- 4. Update the existing FRIDA row (32): order-limit ratio 7:0 -> 5:0 ---
Set-TextValue $ws.Range("H32") "5:0"

# --- 5. Update the grand total (now row 34) to include the new item's price ---
$ws.Range("P34").Value2 = 1468.7000000000001

# --- 6. Refresh the generated timestamp in the footer (now row 35) ---
$ws.Range("A35").Value2 = "Sunday, 21 September, 2025 3:08 PM"
